$d = $word.ActiveDocument
Write-Host "Shapes.Count:" $d.Shapes.Count
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shp = $d.Shapes.Item($i)
    $txt = ""
    if ($shp.TextFrame.HasText) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt.Length -gt 40) {
            $txt = $txt.Substring(0, 40)
        }
    }
    Write-Host $i ": Name=" $shp.Name " Type=" $shp.Type " Text=[" $txt "]"
}
